{"js": "// Merge the two adjacent runs \"`, 1987`\" + \"`)`\" (inside the in-text\n// citation \"(Chiricos, 1987)\") into a single run containing \"`, 1987)`\".\n//\n// The document has the same visible text \"(Chiricos, 1987)\" in two\n// places; only the occurrence inside the paragraph that starts with\n// \"Although most research finds a positive link...\" is still split\n// across two runs, so we scope the search to that paragraph before\n// rewriting the matched range (re-typing identical text over a found\n// range merges the runs it spans into one, which is exactly the\n// run-join the diff shows).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Although most research finds a positive link\") !== -1\n);\n\nif (target) {\n  const hits = target.search(\", 1987)\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(\", 1987)\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Merge the two adjacent runs \"`, 1987`\" + \"`)`\" (inside the in-text\n# citation \"(Chiricos, 1987)\") into a single run containing \"`, 1987)`\".\n#\n# The visible text \"(Chiricos, 1987)\" occurs twice in the document;\n# only the occurrence inside the paragraph that starts with \"Although\n# most research finds a positive link...\" is still split across two\n# runs, so the Find/Replace is scoped to that paragraph's Range only\n# (a document-wide Find would also touch the other, already-merged,\n# occurrence and could even swallow neighbouring runs there).\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*Although most research finds a positive link*\") {\n    $rng = $p.Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = \", 1987)\"\n    $find.Replacement.Text = \", 1987)\"\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    break\n  }\n}\n"}
